# Update the date line at the top of the document.
$d = $word.ActiveDocument
$d.Content.Find.Execute("2025-03-11 Tuesday", $true, $false, $false, $false, $false,
                         $true, 1, $false, "2025-03-12 Wednesday", 2)

# Update the division-problem answers inside the single results table.
# The table has 20 rows (5 data rows interleaved with 3 blank rows each);
# the data rows are 1, 5, 9, 13, 17 (1-based), each with 5 columns.
$t = $d.Tables.Item(1)

$answers = @{
    1  = @("82÷6=13, 4", "44÷3=14, 2", "61÷2=30, 1", "95÷7=13, 4", "97÷3=32, 1")
    5  = @("86÷4=21, 2", "49÷2=24, 1", "78÷2=39, 0", "69÷9=7, 6", "46÷8=5, 6")
    9  = @("53÷5=10, 3", "98÷5=19, 3", "33÷4=8, 1", "22÷5=4, 2", "48÷2=24, 0")
    13 = @("99÷8=12, 3", "76÷3=25, 1", "35÷5=7, 0", "68÷8=8, 4", "12÷5=2, 2")
    17 = @("56÷4=14, 0", "20÷7=2, 6", "32÷3=10, 2", "67÷7=9, 4", "21÷8=2, 5")
}

foreach ($rowIndex in $answers.Keys) {
    $values = $answers[$rowIndex]
    for ($col = 1; $col -le 5; $col++) {
        $cell = $t.Cell($rowIndex, $col)
        $cell.Range.Text = $values[$col - 1]
    }
}
